$wb = $excel.ActiveWorkbook

# --- Update the summary text on sheet "Hoja1" ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.33 = 8878.33 pesos`n✅ 8878.33 pesos = 2.32 = 955.97 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Update the rate figures on sheet "tasas" ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 429.7
$ws2.Range("O10").Value = 3815.02

$ws2.Range("N12").Value = 3831
$ws2.Range("O12").Value = 412.501
